# "Gatinhos para doação" cat-adoption sign-up export — replaces the
# placeholder "Hello World !" sheet with the real donation-sheet dump:
# a merged title row, a green header row, and two zebra-striped data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- whole-sheet base font: Arial 10 (was Calibri 11) -------------------
$ws.Cells.Font.Size = 10
$ws.Cells.Font.Name = "Arial"

# --- clear the old placeholder (content + formatting) ---------------------
$ws.Range("A1").Clear()

# --- title row (merged B3:J3) --------------------------------------------
$ws.Range("B3").Value = "Gatinhos para doação - 2022-07-29 21:43:42"
$ws.Range("B3:J3").Merge()
$ws.Range("B3").Font.Bold = $true
$ws.Range("B3").Font.Size = 12
$ws.Range("B3").Font.Name = "Arial"
$ws.Range("B3").HorizontalAlignment = -4108

# --- header row (B4:J4) ---------------------------------------------------
$ws.Cells.Item(4, 3).Value = "Nome"
$ws.Cells.Item(4, 4).Value = "E-Mail"
$ws.Cells.Item(4, 5).Value = "Telefone"
$ws.Cells.Item(4, 6).Value = "Data De Nascimento"
$ws.Cells.Item(4, 7).Value = "Profissão"
$ws.Cells.Item(4, 8).Value = "Sexo"
$ws.Cells.Item(4, 9).Value = "Cidade"
$ws.Cells.Item(4, 10).Value = "Estado"
$ws.Range("B4:J4").Interior.Color = 11854022
$ws.Range("B4:J4").HorizontalAlignment = -4108

# --- data row 1 (B5:J5) — Pedro Henrique M. Virtuozo -----------------------
$ws.Cells.Item(5, 2).Value = 1
$ws.Cells.Item(5, 3).Value = "Pedro Henrique M. Virtuozo"
$ws.Cells.Item(5, 4).Value = "pedrovirtuozo@gmail.com"
$ws.Cells.Item(5, 5).Value = 5548996234350
$ws.Range("F5").NumberFormat = "@"
$ws.Cells.Item(5, 6).Value = "2004-05-19"
$ws.Cells.Item(5, 7).Value = "Estudante"
$ws.Cells.Item(5, 8).Value = "OUTR"
$ws.Cells.Item(5, 9).Value = "Criciúma"
$ws.Cells.Item(5, 10).Value = "SC"
$ws.Range("B5:J5").Interior.Color = 15917529
$ws.Range("B5:J5").HorizontalAlignment = -4108

# --- data row 2 (B6:J6) — cadu ---------------------------------------------
$ws.Cells.Item(6, 2).Value = 2
$ws.Cells.Item(6, 3).Value = "cadu"
$ws.Cells.Item(6, 4).Value = "cadedu@gmail.com"
$ws.Cells.Item(6, 5).Value = 356457467
$ws.Range("F6").NumberFormat = "@"
$ws.Cells.Item(6, 6).Value = "2022-02-01"
$ws.Cells.Item(6, 7).Value = "Estudante"
$ws.Cells.Item(6, 8).Value = "MASC"
$ws.Cells.Item(6, 9).Value = "Criciúma"
$ws.Cells.Item(6, 10).Value = "SC"
$ws.Range("B6:J6").Interior.Color = 14395790
$ws.Range("B6:J6").HorizontalAlignment = -4108

# --- column widths (character units; COM pads +5px ≈ 0.8333) --------------
$ws.Columns.Item(2).ColumnWidth = 1.1666666666666667   # B = 2
$ws.Columns.Item(3).ColumnWidth = 29.166666666666668   # C = 30
$ws.Columns.Item(4).ColumnWidth = 26.166666666666668   # D = 27
$ws.Columns.Item(5).ColumnWidth = 14.166666666666666   # E = 15
$ws.Columns.Item(6).ColumnWidth = 20.166666666666668   # F = 21
$ws.Columns.Item(7).ColumnWidth = 10.166666666666666   # G = 11
$ws.Columns.Item(8).ColumnWidth = 4.166666666666667    # H = 5
$ws.Columns.Item(9).ColumnWidth = 9.166666666666666    # I = 10
$ws.Columns.Item(10).ColumnWidth = 6.166666666666667   # J = 7

# --- selection lands on the last data row, matching the saved file --------
[void]$ws.Range("B6:J6").Select()
